$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update existing "Actual" values for rows 175-176 and add new rows 177-190
# ---------------------------------------------------------------------------
$actual = @{
  175 = 5534
  176 = 5561
  177 = 5556
  178 = 5590
  179 = 5614
  180 = 5612
  181 = 5635
  182 = 5660
  183 = 5682
  184 = 5746
  185 = 5713
  186 = 5706
  187 = 5743
  188 = 5770
  189 = 5806
  190 = 5858
}

foreach ($r in 175..190) {
  $ws.Cells.Item($r, 3).Value = $actual[$r]
  $ws.Cells.Item($r, 3).NumberFormat = "#,##0"
}

# ---------------------------------------------------------------------------
# 2. Fill the "Daily" column (D) down through row 190
# ---------------------------------------------------------------------------
$ws.Range("D177:D190").Formula = "=C177-C176"
$ws.Range("D177:D190").NumberFormat = "0"

# ---------------------------------------------------------------------------
# 3. Fill the "Average/Week" column (E) down through row 190
# ---------------------------------------------------------------------------
$ws.Range("E177:E190").Formula = "=(C177-C170)/7"
$ws.Range("E177:E190").NumberFormat = "0"

# ---------------------------------------------------------------------------
# 4. Fill the "Average/Ever" column (F) down through row 190
# ---------------------------------------------------------------------------
foreach ($r in 177..190) {
  $formula = "=SUM(`$D`$2:D" + $r + ")/COUNT(`$D`$2:D" + $r + ")"
  $ws.Cells.Item($r, 6).Formula = $formula
  $ws.Cells.Item($r, 6).NumberFormat = "0"
}

# ---------------------------------------------------------------------------
# 5. Extend the Poly-2 (G) and Poly-3 (H) projection formulas through row 238
# ---------------------------------------------------------------------------
$ws.Range("G177:G238").Formula = "=0.2129*B177^2+0.0613*B177+3"
$ws.Range("G177:G238").NumberFormat = "0"

$ws.Range("H177:H238").Formula = "=-0.0008*B177^3+0.3777*B177^2-7.9187*B177+3"
$ws.Range("H177:H238").NumberFormat = "0"

# ---------------------------------------------------------------------------
# 6. Update the saved selection on the worksheet
# ---------------------------------------------------------------------------
$ws.Range("R37").Select()
